$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45309
$ws.Range("D28").Value = 42870
$ws.Range("D29").Value = 57191
